$d = $word.ActiveDocument

$invoiceMarker = [string][char]36 + '{invoice}'
$narrationMarker = [string][char]36 + '{narration}'

# Locate the paragraph that contains the literal dollar-brace "invoice"
# placeholder (sits between the "table" placeholder paragraph and the
# "narration" placeholder paragraph, right before the signature block).
$invoicePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $invoiceMarker) {
        $invoicePara = $p
        break
    }
}

if ($invoicePara -ne $null) {
    # Delete the whole paragraph (its text plus its own paragraph mark), which
    # merges its old position into the following (placeholder "narration")
    # paragraph.
    $invoicePara.Range.Delete()
}

# Find the placeholder "narration" paragraph again post-delete and drop a
# "_GoBack" bookmark right at its start (before the run). Adding a bookmark
# with a name that already exists elsewhere in the document moves it here,
# so this both creates the new bookmark and removes the old one that used to
# sit near the signature line, in a single step.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $narrationMarker) {
        $startRange = $d.Range($p.Range.Start, $p.Range.Start)
        $d.Bookmarks.Add("_GoBack", $startRange)
        break
    }
}
